$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1027
$ws1.Range("F5").Value = 2775
$ws1.Range("F9").Value = 120
$ws1.Range("F11").Value = 67
$ws1.Range("F12").Value = 2588
$ws1.Range("F13").Value = 754

# Sheet "全部类型" (sheet4.xml) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1027
$ws4.Range("F6").Value = 2775
$ws4.Range("F11").Value = 120
$ws4.Range("F13").Value = 67
$ws4.Range("F14").Value = 2588
$ws4.Range("F15").Value = 754
